# "Add files via upload" / "Fixed missing column errors"
#
# Rows 22 and 23 on the active sheet had an erroneous extra "d13C7"
# value duplicated into column C (cells C22/C23). Clear those cells so
# the column lines back up correctly, then leave the selection where the
# author left it when saving (C26).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C22").ClearContents()
$ws.Range("C23").ClearContents()

$ws.Range("C26").Select()
